# Edit script implementing the "Pitch e Documento IOT" commit.
# Applies targeted Find & Replace operations against $word.ActiveDocument
# to reproduce the OOXML diff (text content insertions + proofErr/run cleanup).

$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $findText"
    }
}

# 1) "...onde apos o cadastro do usuario" + proofErr("o mesmo") + "pode criar..."
#    -> merge into a single run (no visible text change, just removes the
#       grammar-check proofErr wrapper around "o mesmo").
Replace-All "cadastro do usuário o mesmo pode criar" "cadastro do usuário o mesmo pode criar"

# 2) "...tamanho, " + proofErr("material, ") + "etc"(spell-checked) + "…" + " "
#    -> remove the grammar proofErr wrapping "material, " (merge with
#       "tamanho, "), but keep the spell-check wrap around "etc" intact.
Replace-All "tamanho, material, " "tamanho, material, "

# ...then merge the trailing ellipsis + space (after "etc") into one run,
#    without touching the "etc" spell-check boundary. Scope the search to
#    start after "Após esta tela" so we hit this specific paragraph.
$anchor = $d.Content
$anchor.Find.Execute("Após esta tela", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterAnchor = $d.Range($anchor.End, $d.Content.End)
$ok = $afterAnchor.Find.Execute("… ", $true, $false, $false, $false, $false, $true, 1, $false, "… ", 2)
if (-not $ok) {
    Write-Output "NOT FOUND: ellipsis after 'Após esta tela'"
}

# 3) "...dividos em Departamentos " + proofErr("e também") + " em " + Tags(spell-checked)
#    -> merge "Departamentos " + "e também" + " em " into a single run,
#       keep the spell-check wrap around "Tags" intact.
Replace-All "Departamentos e também em " "Departamentos e também em "

# 4) Expand the final sentence of the "As cotações funcionarão..." paragraph
#    with the new content about supplier pricing / acceptance flow.
Replace-All "deverá aceitar a cotação, liberando assim os meios de contato entre ambas as partes e dando a possibilidade de avaliações de preço, qualidade e entrega da parte do comprador para com o fornecedor." "deverá aceitar a cotação porém o mesmo não saberá até este momento o valor que o cliente está disposto a pagar pelo produto. Sendo assim, o fornecedor deverá aceitar ou recusar cotações de acordo com quantidades, prazos, produtos, etc.. Após aprovada a cotação, o fornecedor deverá informar o preço unitário que ele está disposto a vender o produto e caso este seja menor ou igual ao valor solicitado pelo comprador, ou dentro de algum range específico de preço (por exemplo uma margem de até 10% acima), o orçamento será enviado para o comprador, liberando assim os meios de contato entre ambas as partes e dando a possibilidade de avaliações de preço, qualidade e entrega da parte do comprador para com o fornecedor."

# 5) "...quantidade indicada " + proofErr("pelo mesmo") + ". "
#    -> merge into "quantidade indicada pelo mesmo. "
Replace-All "quantidade indicada pelo mesmo" "quantidade indicada pelo mesmo"

# 6) Insert " e informar o seu preço" after "aquelas condições" in the
#    "A nova ideia consiste..." paragraph.
Replace-All "aquelas condições, garantindo" "aquelas condições e informar o seu preço, garantindo"

# 7) "...de prazo, qualidade, " + proofErr("avaliações, ") + "etc"(spell-checked) + "…"
#    -> merge "qualidade, " + "avaliações, " into a single run (drop the
#       grammar proofErr wrap), and also drop the trailing gramEnd marker
#       that followed the ellipsis (no text change there).
Replace-All "de prazo, qualidade, avaliações, " "de prazo, qualidade, avaliações, "

# 8) NUMPY paragraph: "...valores " + proofErr("totais, ") + "etc"(spell-checked) + "…" + ";"
#    -> merge "valores " + "totais, " into one run, then separately merge
#       the trailing ellipsis with the following ";" run.
Replace-All "valores totais, " "valores totais, "
Replace-All "etc…;" "etc…;"

# 9) "...para que " + proofErr("a mesma") + " atenda a sua necessidade..."
#    -> merge into "para que a mesma atenda"
Replace-All "para que a mesma atenda" "para que a mesma atenda"
